# Adds a new oracle row (row 7) to Sheet1 of the oracles-dataset workbook.
# The new row documents PolynomialFunction.differentiate(double[]) and is
# modelled on the existing rows for the same class (rows 3/4 for most
# columns' styling/content, rows 5/6 for the javadocTag quote-prefix style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

# --- New text introduced by this edit -------------------------------------

$oracle = "true ? Arrays.stream(coefficients).noneMatch(jdVar -> jdVar == null) : true;"

$javadocTag = "@return the coefficients of the derivative or {@code null} if coefficients has length 1."

$methodJavadoc = @"
/**
     * Returns the coefficients of the derivative of the polynomial with the given coefficients.
     *
     * @param coefficients Coefficients of the polynomial to differentiate.
     * @return the coefficients of the derivative or {@code null} if coefficients has length 1.
     * @throws NoDataException if {@code coefficients} is empty.
     * @throws NullArgumentException if {@code coefficients} is {@code null}.
     */
"@

$methodSourceCode = @"
protected static double[] differentiate(double[] coefficients)
        throws NullArgumentException, NoDataException {
        MathUtils.checkNotNull(coefficients);
        int n = coefficients.length;
        if (n == 0) {
            throw new NoDataException(LocalizedFormats.EMPTY_POLYNOMIALS_COEFFICIENTS_ARRAY);
        }
        if (n == 1) {
            return new double[]{0};
        }
        double[] result = new double[n - 1];
        for (int i = n - 1; i > 0; i--) {
            result[i - 1] = i * coefficients[i];
        }
        return result;
    }
"@

$methodJavadocValues = "[1; int]"
$methodArguments = "[coefficients;;double\[\]]"

# --- Columns whose values already exist elsewhere in the sheet: reuse them
# via .Text so the workbook's shared-string table is deduplicated exactly
# like a native Excel edit would do. ---------------------------------------

$oracleType          = $ws.Range("C5").Text
$projectName         = $ws.Range("D3").Text
$packageName         = $ws.Range("E3").Text
$className           = $ws.Range("F3").Text
$classJavadoc        = $ws.Range("J3").Text
$classSourceCode     = $ws.Range("K3").Text
$tokensGeneralGrammar                                     = $ws.Range("L2").Text
$tokensGeneralValuesGlobalDictionary                      = $ws.Range("M2").Text
$tokensProjectClasses                                     = $ws.Range("N3").Text
$tokensProjectClassesNonPrivateStaticNonVoidMethods       = $ws.Range("O2").Text
$tokensProjectClassesNonPrivateStaticAttributes           = $ws.Range("P2").Text
$tokensMethodVariablesNonPrivateNonStaticNonVoidMethods   = $ws.Range("S3").Text

# --- Write the new row ------------------------------------------------------

$ws.Cells.Item($row, 1).Value = 6                     # A7 id
$ws.Cells.Item($row, 2).Value = $oracle                # B7 oracle
$ws.Cells.Item($row, 3).Value = $oracleType            # C7 oracleType
$ws.Cells.Item($row, 4).Value = $projectName           # D7 projectName
$ws.Cells.Item($row, 5).Value = $packageName           # E7 packageName
$ws.Cells.Item($row, 6).Value = $className             # F7 className
$ws.Cells.Item($row, 7).Value = $javadocTag            # G7 javadocTag
$ws.Cells.Item($row, 8).Value = $methodJavadoc         # H7 methodJavadoc
$ws.Cells.Item($row, 9).Value = $methodSourceCode      # I7 methodSourceCode
$ws.Cells.Item($row, 10).Value = $classJavadoc         # J7 classJavadoc
$ws.Cells.Item($row, 11).Value = $classSourceCode      # K7 classSourceCode
$ws.Cells.Item($row, 12).Value = $tokensGeneralGrammar # L7 tokensGeneralGrammar
$ws.Cells.Item($row, 13).Value = $tokensGeneralValuesGlobalDictionary # M7
$ws.Cells.Item($row, 14).Value = $tokensProjectClasses # N7 tokensProjectClasses
$ws.Cells.Item($row, 15).Value = $tokensProjectClassesNonPrivateStaticNonVoidMethods # O7
$ws.Cells.Item($row, 16).Value = $tokensProjectClassesNonPrivateStaticAttributes     # P7
$ws.Cells.Item($row, 17).Value = $methodJavadocValues  # Q7 tokensMethodJavadocValues
$ws.Cells.Item($row, 18).Value = $methodArguments      # R7 tokensMethodArguments
$ws.Cells.Item($row, 19).Value = $tokensMethodVariablesNonPrivateNonStaticNonVoidMethods # S7

# U7 / V7 stay empty but pick up the wrap-text style below.

# --- Formatting --------------------------------------------------------------

# G7 uses the "quote prefix" style (like G5/G6, whose javadocTag also starts
# with '@'); copy that format only (not the value) from G6.
$ws.Range("G6").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null

# H7, I7, J7, K7, N7, R7, S7, U7, V7 use the wrap-text style (like the rest
# of the sheet's long free-text columns). U7/V7 stay empty but still carry
# the style, mirroring U4/V4 above.
foreach ($col in @("H", "I", "J", "K", "N", "R", "S", "U", "V")) {
    $ws.Range($col + "7").WrapText = $true
}

$ws.Application.CutCopyMode = $false

# Row height matches the other long-content rows in the sheet.
$ws.Rows.Item($row).RowHeight = 409.6

# Selection matches the position left behind by this edit.
$ws.Range("S6").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 19
$ws.Range("U12").Select() | Out-Null
